# Updated the input files in sixteen_tests to have strain_log2_expression instead of
# just strain and then ran the files and saved the outputs in sixteen_tests_output.
#
# This script:
#  1. Renames worksheet "wt"    -> "wt_log2_expression"
#  2. Renames worksheet "dcin5" -> "dcin5_log2_expression"
#  3. Moves the selected cell on the dcin5 sheet from O9 to E38
#     (while preserving which sheet/tab is active in the workbook)

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it after
# touching the "dcin5" sheet's selection (selecting a range on a sheet
# activates that sheet as a side effect).
$originalActiveSheetName = $wb.ActiveSheet.Name

# --- Rename sheets ---------------------------------------------------
$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"

# --- Update the selection on the renamed "dcin5_log2_expression" sheet
$dcin5 = $wb.Worksheets.Item("dcin5_log2_expression")
$dcin5.Range("E38").Select() | Out-Null

# --- Restore the originally active sheet/tab -------------------------
$wb.Worksheets.Item($originalActiveSheetName).Activate() | Out-Null

Write-Host "Sheets renamed and selection updated."
